# Removed Encounter.header data, keep order in mapped tables
#
# The EHDSEncounter.header.* sub-items (rows 4-14 in column A) are removed
# from the mapping sheet. Since column B ("zib"/Encounter mapped values) had
# no data in those rows, deleting the entire rows shifts everything below
# upward by 11 rows while preserving the relative order/association of the
# remaining A/B pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete rows 4 through 14 (EHDSEncounter.header.subject .. EHDSEncounter.presentedForm),
# shifting the remaining rows up.
$ws.Range("A4:B14").EntireRow.Delete()

# Restore the selection to match the post-edit cursor position.
$ws.Range("B9").Select() | Out-Null
